# Update the "Broken" sheet (broken_links report) and the "Orphan" sheet
# (orphan_pages report) for hosninsurance.ae with the new social-media /
# internal-link data, and trim each sheet to its new (smaller) row count.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Broken  (Source Page URL | Broken Link URL | Anchor Text /
#                 Current Value | Link Type | Status Code)
# ---------------------------------------------------------------------
$broken = $wb.Worksheets.Item("Broken")

$brokenRows = @(
    @{ A = "https://hosninsurance.ae/";                                       B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/services-hosn-al-sharjah-insurance";      B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/about-hosn-al-sharjah-insurance";         B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/contact-us";                             B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/blogs";                                  B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/services-hosn-al-sharjah-insurance/";    B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/about-hosn-al-sharjah-insurance/";       B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/blog-article";                          B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/author/0xdanielimad";                   B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" },
    @{ A = "https://hosninsurance.ae/category/uncategorized";                B = "https://www.instagram.com/hosninsurance"; C = "(no anchor text)"; D = "External"; E = "429" }
)

$r = 2
foreach ($row in $brokenRows) {
    $broken.Cells.Item($r, 1).Value = $row.A
    $broken.Cells.Item($r, 2).Value = $row.B
    $broken.Cells.Item($r, 3).Value = $row.C
    $broken.Cells.Item($r, 4).Value = $row.D
    $broken.Cells.Item($r, 5).Value = $row.E
    $r = $r + 1
}

# The original sheet had data through row 26; the updated report only has
# data through row 11, so remove the now-unused trailing rows.
$broken.Range("A12:E26").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------
# Sheet: Orphan  (Orphan Page URL | Found in Sitemap? | Internally Linked?)
# ---------------------------------------------------------------------
$orphan = $wb.Worksheets.Item("Orphan")

$orphanRows = @(
    @{ A = "https://hosninsurance.ae/contact-us";                          B = "Yes"; C = "Yes" },
    @{ A = "https://hosninsurance.ae/about-hosn-al-sharjah-insurance";      B = "Yes"; C = "Yes" },
    @{ A = "https://hosninsurance.ae/category/uncategorized";               B = "Yes"; C = "Yes" },
    @{ A = "https://hosninsurance.ae";                                      B = "Yes"; C = "No"  },
    @{ A = "https://hosninsurance.ae/blog-article";                        B = "Yes"; C = "Yes" },
    @{ A = "https://hosninsurance.ae/blogs";                                B = "Yes"; C = "Yes" },
    @{ A = "https://hosninsurance.ae/services-hosn-al-sharjah-insurance";   B = "Yes"; C = "Yes" }
)

$r = 2
foreach ($row in $orphanRows) {
    $orphan.Cells.Item($r, 1).Value = $row.A
    $orphan.Cells.Item($r, 2).Value = $row.B
    $orphan.Cells.Item($r, 3).Value = $row.C
    $r = $r + 1
}

# The original sheet had data through row 16; the updated report only has
# data through row 8, so remove the now-unused trailing rows.
$orphan.Range("A9:C16").EntireRow.Delete() | Out-Null
